$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at G (shifts the existing "s1cDNAProtocol" column from G to H)
$ws.Range("G1").EntireColumn.Insert()

# New header cell G1: "polyAIsolationProtocol"
$ws.Range("G1").Style = "Normal"
$ws.Range("G1").Value = "polyAIsolationProtocol"
$ws.Range("G1").Font.Name = "Helvetica"
$ws.Range("G1").Font.Color = 3025188

# New column data G2:G29: "catcher"
$ws.Range("G2:G29").Value = "catcher"

# Column widths (author resized C:H after adding the new column; values chosen
# so the engine's internal 1/6-character quantization lands on the closest
# achievable width to the target stored width)
$ws.Columns.Item(3).ColumnWidth = 14.833333333333334
$ws.Columns.Item(4).ColumnWidth = 10
$ws.Columns.Item(5).ColumnWidth = 13.166666666666666
$ws.Columns.Item(6).ColumnWidth = 18.5
$ws.Columns.Item(7).ColumnWidth = 18.5
$ws.Columns.Item(8).ColumnWidth = 12.833333333333334

# Restore selection to where the user was working
$ws.Range("G3:G29").Select()
